# "back to original function"
# Reverts the PSNR worksheet from the 4-column comparison (LSB / LSB-DES /
# LSB-pair / LSB-pair1) back to the original 3-column layout (LSB / LSB-pair-ultar),
# restoring column D's values from column B and dropping column E entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop column E (LSB-pair1) completely - data shifts left, dimension becomes A1:D101.
$ws.Columns.Item(5).Delete()

# Column C / D headers both become "LSB-pair-ultar" (C1 was "LSB-DES",
# D1 was "LSB-pair" pointing at a now-removed string).
$ws.Range("C1").Value = "LSB-pair-ultar"
$ws.Range("D1").Value = "LSB-pair-ultar"

# Column C gets refreshed data values for the first 10 data rows.
$ws.Range("C2").Value2 = 44554
$ws.Range("C3").Value2 = 43921
$ws.Range("C4").Value2 = 44599
$ws.Range("C5").Value2 = 44710
$ws.Range("C6").Value2 = 44016
$ws.Range("C7").Value2 = 43010
$ws.Range("C8").Value2 = 37354
$ws.Range("C9").Value2 = 43426
$ws.Range("C10").Value2 = 43283
$ws.Range("C11").Value2 = 45620

# Column D (last remaining column, formerly housing LSB-pair data) now mirrors
# column B for those same rows.
$ws.Range("D2").Value2 = $ws.Range("B2").Value2
$ws.Range("D3").Value2 = $ws.Range("B3").Value2
$ws.Range("D4").Value2 = $ws.Range("B4").Value2
$ws.Range("D5").Value2 = $ws.Range("B5").Value2
$ws.Range("D6").Value2 = $ws.Range("B6").Value2
$ws.Range("D7").Value2 = $ws.Range("B7").Value2
$ws.Range("D8").Value2 = $ws.Range("B8").Value2
$ws.Range("D9").Value2 = $ws.Range("B9").Value2
$ws.Range("D10").Value2 = $ws.Range("B10").Value2
$ws.Range("D11").Value2 = $ws.Range("B11").Value2

# Restore the original selection rectangle now that only B:D hold data.
$ws.Range("B2:C101").Select()
